$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of rows 21 and 22 for columns A, B, E, F, G, H, Q, R
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr21 = "${col}21"
    $addr22 = "${col}22"
    $val21 = $ws.Range($addr21).Value2
    $val22 = $ws.Range($addr22).Value2
    $ws.Range($addr21).Value2 = $val22
    $ws.Range($addr22).Value2 = $val21
}
